{"js": "// Update the date heading paragraph (first paragraph of the body).\nconst body = context.document.body;\nconst titleHits = body.search(\"2024-02-01 Thursday\", { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\n\nif (titleHits.items.length > 0) {\n  titleHits.items[0].insertText(\"2024-02-02 Friday\", \"Replace\");\n}\n\n// Update the practice-problem table: every 4th row (0, 4, 8, 12, 16) holds\n// five \"NN\u00f7N=\" cells; the rest are intentionally blank spacer rows. Each of\n// those 25 cells is replaced in place (by row/column position, since some\n// values repeat with different targets) with the new problem text.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = {\n  0: [\"16\u00f76=\", \"11\u00f75=\", \"94\u00f72=\", \"30\u00f72=\", \"57\u00f77=\"],\n  4: [\"12\u00f76=\", \"79\u00f72=\", \"93\u00f72=\", \"53\u00f78=\", \"30\u00f72=\"],\n  8: [\"54\u00f75=\", \"59\u00f79=\", \"57\u00f79=\", \"36\u00f78=\", \"77\u00f79=\"],\n  12: [\"75\u00f72=\", \"48\u00f75=\", \"81\u00f77=\", \"83\u00f76=\", \"93\u00f79=\"],\n  16: [\"49\u00f79=\", \"83\u00f74=\", \"30\u00f74=\", \"46\u00f77=\", \"65\u00f77=\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const r = parseInt(rowIndex, 10);\n  const rowValues = newValues[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = rowValues[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph of the document).\n$find = $d.Content.Find\n$find.Execute(\"2024-02-01 Thursday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2024-02-02 Friday\", 2) | Out-Null\n\n# Update the practice-problem table: every 4th row (1,5,9,13,17 in 1-based\n# COM indexing) holds five \"NN\u00f7N=\" cells; the other rows are intentionally\n# blank spacer rows. Replace each of those 25 cells in place (by row/column\n# position, since a couple of values repeat with different targets).\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"16\u00f76=\", \"11\u00f75=\", \"94\u00f72=\", \"30\u00f72=\", \"57\u00f77=\")\n    5  = @(\"12\u00f76=\", \"79\u00f72=\", \"93\u00f72=\", \"53\u00f78=\", \"30\u00f72=\")\n    9  = @(\"54\u00f75=\", \"59\u00f79=\", \"57\u00f79=\", \"36\u00f78=\", \"77\u00f79=\")\n    13 = @(\"75\u00f72=\", \"48\u00f75=\", \"81\u00f77=\", \"83\u00f76=\", \"93\u00f79=\")\n    17 = @(\"49\u00f79=\", \"83\u00f74=\", \"30\u00f74=\", \"46\u00f77=\", \"65\u00f77=\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowValues = $newValues[$rowIndex]\n    for ($c = 1; $c -le $rowValues.Length; $c++) {\n        $t.Cell($rowIndex, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
